$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'running tights women compression',
    'running tights women pocket',
    'running tights women pockets',
    'running tights womens',
    'running tights zipper',
    'running touch waist pack',
    'running tracks',
    'running training equipment',
    'running training gear',
    'running underwear for women',
    'running vest cold weather',
    'running waist',
    'running waist pack',
    'running waistband iphone 6',
    'running wear',
    'running wear for women',
    'running wear women',
    'running wearable',
    'running wind pants women',
    'running winter gear men',
    'running women clothing',
    'running women cold weather',
    'running womens clothes',
    'running workout pants',
    'running yoga pants',
    'running yoga pants for women',
    'runnning pants',
    'runnung pants for women',
    'rustic basketball',
    'rustic basketball hoop',
    's a n tight',
    's clothes',
    'saddle up game',
    'saddleseat apparel',
    'saddleseat jodhpurs',
    'saddleseat suit',
    'safe place shirt',
    'safe sport knee brace',
    'sagging pants',
    'saggy pants for women',
    'sartorius muscle',
    'sausage run game',
    'sausage skin',
    'schooling tights',
    'schooling tights equestrian womens',
    'sciatic ball',
    'sciatic nerve ball',
    'sciatica band',
    'sciatica compression',
    'scrunch pants tan',
    'seam guide tape',
    'seam panty hose',
    'seam tights for women',
    'seamless capri leggings',
    'seamless compression leggings',
    'seamless compression leggings for women',
    'second skin clothing women',
    'second skin compression men',
    'second skin compression pants men',
    'second skin compression women',
    'second skin leggings women',
    'second skin long underwear',
    'second skin mens compression',
    'second skin pants men',
    'second skin underwear men',
    'second skin womens pants',
    'second skin workout clothes women',
    'secret leggings',
    'secretly yours',
    'see and wear',
    'see line woman',
    'see through leggings low rise',
    'see thru yoga pants',
    'seem tights',
    'self protection devices',
    'self protection devices for women',
    'self support',
    'send a hug',
    'sensitive skin kinesiology tape',
    'shade for soccer games',
    'shape tanks for women',
    'shape up leggings',
    'shape up tape',
    'shapermint for women leggings',
    'shapermint leggings',
    'shaping capri leggings',
    'sheer black knee highs',
    'sheer capri leggings for women',
    'sheer clothing for women',
    'sheer compression knee highs',
    'sheer dark green tights',
    'sheer genius',
    'sheer leggings',
    'sheer leggings for women',
    'sheer leggings women',
    'sheer nylon leggings',
    'sheer panties xl',
    'sheer red tights',
    'sheer runner',
    'sheer shorts for women'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
